$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet so stale cells/rows from the old layout do not
# linger (the new layout re-flows several rows).
$ws.Cells.Clear()

# ---- Section 1: Clock / prescaler -------------------------------------
$ws.Range("A1").Value2 = 4
$ws.Range("A1").Interior.Color = 65535
$ws.Range("B1").Value2 = "MHz"
$ws.Range("C1").Value2 = "Clock"

$ws.Range("A2").Value2 = 1024
$ws.Range("A2").Interior.Color = 65535
$ws.Range("C2").Value2 = "Prescaler"

$ws.Range("A3").Formula = "=A1*1000/A2"
$ws.Range("A3").NumberFormat = "0.00"
$ws.Range("B3").Value2 = "kHz"
$ws.Range("C3").Value2 = "Prescaled clock"

$ws.Range("A4").Formula = "=1/(A3)"
$ws.Range("B4").Value2 = "ms"

# ---- Section 2: interrupt timing ---------------------------------------
$ws.Range("A6").Value2 = "interrupt"

$ws.Range("A7").Value2 = 750
$ws.Range("A7").Interior.Color = 65535
$ws.Range("B7").Value2 = "ms"
$ws.Range("C7").Value2 = 256
$ws.Range("C7").Interior.Color = 65535
$ws.Range("D7").Value2 = "iterations"

$ws.Range("A8").Formula = "=A7*A1*1000000/A2/1000"
$ws.Range("B8").Value2 = "iterations"
$ws.Range("C8").Formula = "=C7*A4"
$ws.Range("D8").Value2 = "ms"

# ---- Section 3: TWI / clock prescaler -----------------------------------
$ws.Range("A10").Value2 = "TWI"

$ws.Range("A11").Value2 = "F_CLK"
$ws.Range("B11").Formula = "=A1"
$ws.Range("B11").Interior.ColorIndex = 0
$ws.Range("C11").Value2 = "MHz"

$ws.Range("A12").Value2 = "TWBR"
$ws.Range("B12").Formula = "=(B11*2/(B14/1000)-16)/2"

$ws.Range("A13").Value2 = "TWPS"
$ws.Range("B13").Value2 = 1

$ws.Range("A14").Value2 = "f"
$ws.Range("B14").Value2 = 100
$ws.Range("B14").NumberFormat = "0"
$ws.Range("B14").Interior.Color = 65535
$ws.Range("C14").Value2 = "kHz"

# ---- Misc formatting left over from the original sheet ------------------
$ws.Range("B21").NumberFormat = "0"

$ws.Columns("A:E").ColumnWidth = 9.75

$ws.PageSetup.Orientation = 1

[void]$ws.Range("A11").Select()
